$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 160.334686357684
$ws.Cells.Item(1, 2).Value = 20.6578374005192
$ws.Cells.Item(1, 3).Value = 125.304815138366
$ws.Cells.Item(1, 4).Value = 7.42962435233855
$ws.Cells.Item(1, 5).Value = 141.576353060816
$ws.Cells.Item(1, 6).Value = 72.0289083533124
$ws.Cells.Item(1, 7).Value = 99.7806833590291
$ws.Cells.Item(1, 8).Value = 10.6914648835973
$ws.Cells.Item(1, 9).Value = 29.0826002271299
$ws.Cells.Item(1, 10).Value = 186.48344063502

$ws.Cells.Item(2, 1).Value = 49.5330577015565
$ws.Cells.Item(2, 2).Value = 21.2742932239893
$ws.Cells.Item(2, 3).Value = 42.4062806379079
$ws.Cells.Item(2, 4).Value = 39.5965671351164
$ws.Cells.Item(2, 5).Value = 103.93235157427
$ws.Cells.Item(2, 6).Value = 139.959372458961
$ws.Cells.Item(2, 7).Value = 137.812589638779
$ws.Cells.Item(2, 8).Value = 6.90119434469435
$ws.Cells.Item(2, 9).Value = 139.518385166078
$ws.Cells.Item(2, 10).Value = 165.918658844157

$ws.Cells.Item(3, 1).Value = 180.165873831215
$ws.Cells.Item(3, 2).Value = 131.93953006153
$ws.Cells.Item(3, 3).Value = 186.906220757824
$ws.Cells.Item(3, 4).Value = 58.8757023489455
$ws.Cells.Item(3, 5).Value = 62.2783218800455
$ws.Cells.Item(3, 6).Value = 128.016692087062
$ws.Cells.Item(3, 7).Value = 149.201590264776
$ws.Cells.Item(3, 8).Value = 45.5876905683371
$ws.Cells.Item(3, 9).Value = 118.602181653773
$ws.Cells.Item(3, 10).Value = 62.258485454255

$ws.Cells.Item(4, 1).Value = 35.4854357594091
$ws.Cells.Item(4, 2).Value = 101.26448846481
$ws.Cells.Item(4, 3).Value = 47.5028627773295
$ws.Cells.Item(4, 4).Value = 92.8036313004809
$ws.Cells.Item(4, 5).Value = 0.8067236285688
$ws.Cells.Item(4, 6).Value = 101.407469763145
$ws.Cells.Item(4, 7).Value = 74.5363784369716
$ws.Cells.Item(4, 8).Value = 181.77221239627
$ws.Cells.Item(4, 9).Value = 191.688250653301
$ws.Cells.Item(4, 10).Value = 90.7576557671454

$ws.Cells.Item(5, 1).Value = 94.9912441405427
$ws.Cells.Item(5, 2).Value = 3.69410598822595
$ws.Cells.Item(5, 3).Value = 4.29676138064673
$ws.Cells.Item(5, 4).Value = 190.969530209419
$ws.Cells.Item(5, 5).Value = 126.354855823496
$ws.Cells.Item(5, 6).Value = 29.8405539383369
$ws.Cells.Item(5, 7).Value = 101.439659065306
$ws.Cells.Item(5, 8).Value = 116.551897729073
$ws.Cells.Item(5, 9).Value = 105.142685726817
$ws.Cells.Item(5, 10).Value = 199.370667151814

$ws.Cells.Item(6, 1).Value = 44.2154777442177
$ws.Cells.Item(6, 2).Value = 160.873748250713
$ws.Cells.Item(6, 3).Value = 112.889927585092
$ws.Cells.Item(6, 4).Value = 124.503278231483
$ws.Cells.Item(6, 5).Value = 163.190609478946
$ws.Cells.Item(6, 6).Value = 28.3951562961541
$ws.Cells.Item(6, 7).Value = 33.7516166426947
$ws.Cells.Item(6, 8).Value = 66.4291127894209
$ws.Cells.Item(6, 9).Value = 145.151302472293
$ws.Cells.Item(6, 10).Value = 13.5596609737536

$ws.Cells.Item(7, 1).Value = 122.827318088537
$ws.Cells.Item(7, 2).Value = 54.192992790692
$ws.Cells.Item(7, 3).Value = 92.0892832298247
$ws.Cells.Item(7, 4).Value = 166.824114772875
$ws.Cells.Item(7, 5).Value = 150.998004875611
$ws.Cells.Item(7, 6).Value = 148.268569236746
$ws.Cells.Item(7, 7).Value = 173.77143044666
$ws.Cells.Item(7, 8).Value = 149.602649337427
$ws.Cells.Item(7, 9).Value = 38.7898435065475
$ws.Cells.Item(7, 10).Value = 2.52488181112561

$ws.Cells.Item(8, 1).Value = 65.4229940219889
$ws.Cells.Item(8, 2).Value = 156.040377242509
$ws.Cells.Item(8, 3).Value = 15.2129436913938
$ws.Cells.Item(8, 4).Value = 48.7607293989326
$ws.Cells.Item(8, 5).Value = 70.9274147036147
$ws.Cells.Item(8, 6).Value = 176.471767842989
$ws.Cells.Item(8, 7).Value = 127.642768680883
$ws.Cells.Item(8, 8).Value = 195.936690548405
$ws.Cells.Item(8, 9).Value = 132.52084652545
$ws.Cells.Item(8, 10).Value = 32.4377679417086

$ws.Cells.Item(9, 1).Value = 26.577033021756
$ws.Cells.Item(9, 2).Value = 32.6496925357029
$ws.Cells.Item(9, 3).Value = 140.44500484152
$ws.Cells.Item(9, 4).Value = 119.231514501959
$ws.Cells.Item(9, 5).Value = 18.0430077100373
$ws.Cells.Item(9, 6).Value = 74.6116875086965
$ws.Cells.Item(9, 7).Value = 188.374560879718
$ws.Cells.Item(9, 8).Value = 122.999584545847
$ws.Cells.Item(9, 9).Value = 129.613021821535
$ws.Cells.Item(9, 10).Value = 172.411567332415

$ws.Cells.Item(10, 1).Value = 67.6558531204499
$ws.Cells.Item(10, 2).Value = 8.1072656475507
$ws.Cells.Item(10, 3).Value = 36.6209099239767
$ws.Cells.Item(10, 4).Value = 178.128589679547
$ws.Cells.Item(10, 5).Value = 167.930337678609
$ws.Cells.Item(10, 6).Value = 40.7982513498507
$ws.Cells.Item(10, 7).Value = 111.604822758401
$ws.Cells.Item(10, 8).Value = 37.4726466077718
$ws.Cells.Item(10, 9).Value = 39.9715253338085
$ws.Cells.Item(10, 10).Value = 178.08628658675

$ws.Cells.Item(11, 1).Value = 56.069123491677
$ws.Cells.Item(11, 2).Value = 151.837009727879
$ws.Cells.Item(11, 3).Value = 77.7620542225251
$ws.Cells.Item(11, 4).Value = 102.617803915691
$ws.Cells.Item(11, 5).Value = 133.947673129825
$ws.Cells.Item(11, 6).Value = 88.1751005017083
$ws.Cells.Item(11, 7).Value = 145.660804559319
$ws.Cells.Item(11, 8).Value = 64.1291981861597
$ws.Cells.Item(11, 9).Value = 53.5758635278679
$ws.Cells.Item(11, 10).Value = 186.718841635957

$ws.Cells.Item(12, 1).Value = 100.752387615271
$ws.Cells.Item(12, 2).Value = 37.8149260942894
$ws.Cells.Item(12, 3).Value = 133.908266263971
$ws.Cells.Item(12, 4).Value = 112.713534530584
$ws.Cells.Item(12, 5).Value = 186.982627951998
$ws.Cells.Item(12, 6).Value = 90.1776255528338
$ws.Cells.Item(12, 7).Value = 113.747987949172
$ws.Cells.Item(12, 8).Value = 172.857768727866
$ws.Cells.Item(12, 9).Value = 148.781107062838
$ws.Cells.Item(12, 10).Value = 76.386317366914

$ws.Cells.Item(13, 1).Value = 159.894008357028
$ws.Cells.Item(13, 2).Value = 50.771845900813
$ws.Cells.Item(13, 3).Value = 19.9896275158923
$ws.Cells.Item(13, 4).Value = 66.3782761741329
$ws.Cells.Item(13, 5).Value = 134.869028690676
$ws.Cells.Item(13, 6).Value = 57.3157283744382
$ws.Cells.Item(13, 7).Value = 119.419467318533
$ws.Cells.Item(13, 8).Value = 37.0843540118469
$ws.Cells.Item(13, 9).Value = 80.830391720324
$ws.Cells.Item(13, 10).Value = 30.1291633537641

$ws.Cells.Item(14, 1).Value = 64.8669450845881
$ws.Cells.Item(14, 2).Value = 90.1701220731112
$ws.Cells.Item(14, 3).Value = 155.965165214597
$ws.Cells.Item(14, 4).Value = 154.4345599387
$ws.Cells.Item(14, 5).Value = 176.368644450032
$ws.Cells.Item(14, 6).Value = 74.740023293877
$ws.Cells.Item(14, 7).Value = 154.887638313178
$ws.Cells.Item(14, 8).Value = 37.8272009258285
$ws.Cells.Item(14, 9).Value = 185.283841372134
$ws.Cells.Item(14, 10).Value = 129.867907208329

$ws.Cells.Item(15, 1).Value = 128.950882949378
$ws.Cells.Item(15, 2).Value = 124.245362693558
$ws.Cells.Item(15, 3).Value = 69.423721017979
$ws.Cells.Item(15, 4).Value = 142.894180185578
$ws.Cells.Item(15, 5).Value = 71.6591797171436
$ws.Cells.Item(15, 6).Value = 29.8409270261605
$ws.Cells.Item(15, 7).Value = 74.1989993835795
$ws.Cells.Item(15, 8).Value = 123.907375393392
$ws.Cells.Item(15, 9).Value = 191.145961727549
$ws.Cells.Item(15, 10).Value = 77.7527121257748

$ws.Cells.Item(16, 1).Value = 127.050263400679
$ws.Cells.Item(16, 2).Value = 138.747054030535
$ws.Cells.Item(16, 3).Value = 88.6915395449342
$ws.Cells.Item(16, 4).Value = 163.585207966894
$ws.Cells.Item(16, 5).Value = 18.192278229721
$ws.Cells.Item(16, 6).Value = 5.297277590864
$ws.Cells.Item(16, 7).Value = 131.847382211987
$ws.Cells.Item(16, 8).Value = 11.3837780483923
$ws.Cells.Item(16, 9).Value = 167.748775225016
$ws.Cells.Item(16, 10).Value = 76.6319447553865

$ws.Cells.Item(17, 1).Value = 168.755633183176
$ws.Cells.Item(17, 2).Value = 108.576450547472
$ws.Cells.Item(17, 3).Value = 183.298806465836
$ws.Cells.Item(17, 4).Value = 23.4467001741038
$ws.Cells.Item(17, 5).Value = 121.851896551369
$ws.Cells.Item(17, 6).Value = 10.5822655421599
$ws.Cells.Item(17, 7).Value = 81.8497608796925
$ws.Cells.Item(17, 8).Value = 179.473706325271
$ws.Cells.Item(17, 9).Value = 136.344890080553
$ws.Cells.Item(17, 10).Value = 112.242604658121

$ws.Cells.Item(18, 1).Value = 135.289987239656
$ws.Cells.Item(18, 2).Value = 75.9207870233435
$ws.Cells.Item(18, 3).Value = 187.573927355732
$ws.Cells.Item(18, 4).Value = 18.9131998545086
$ws.Cells.Item(18, 5).Value = 147.435434417536
$ws.Cells.Item(18, 6).Value = 35.6486456634703
$ws.Cells.Item(18, 7).Value = 181.348124882834
$ws.Cells.Item(18, 8).Value = 77.0954473303144
$ws.Cells.Item(18, 9).Value = 194.719096456989
$ws.Cells.Item(18, 10).Value = 105.028101664515

$ws.Cells.Item(19, 1).Value = 183.116728990859
$ws.Cells.Item(19, 2).Value = 195.512091925141
$ws.Cells.Item(19, 3).Value = 45.9383922842976
$ws.Cells.Item(19, 4).Value = 3.0776795945492
$ws.Cells.Item(19, 5).Value = 103.078899953085
$ws.Cells.Item(19, 6).Value = 126.119891054053
$ws.Cells.Item(19, 7).Value = 1.47858252817699
$ws.Cells.Item(19, 8).Value = 192.379957247702
$ws.Cells.Item(19, 9).Value = 136.242281708979
$ws.Cells.Item(19, 10).Value = 171.071366859168

$ws.Cells.Item(20, 1).Value = 142.89264108189
$ws.Cells.Item(20, 2).Value = 143.503860264786
$ws.Cells.Item(20, 3).Value = 70.0784257008128
$ws.Cells.Item(20, 4).Value = 108.651896616748
$ws.Cells.Item(20, 5).Value = 161.112274025153
$ws.Cells.Item(20, 6).Value = 20.374432401906
$ws.Cells.Item(20, 7).Value = 140.946556227722
$ws.Cells.Item(20, 8).Value = 45.9770208438751
$ws.Cells.Item(20, 9).Value = 21.0422836342092
$ws.Cells.Item(20, 10).Value = 61.0769141749837

